# Apply updates to the "metricas_analise_estoque" sheet:
#  - rename several header labels in row 1
#  - update the analysis timestamp and several computed metrics in row 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header label renames -----------------------------------------
$ws.Range("C1").Value = "TOTAL SKU COM HISTORICO > 1 ANO"
$ws.Range("D1").Value = "%SKU COM COM HISTORICO > 1 ANO"
$ws.Range("E1").Value = "TOTAL SKU COM HISTORICO < 1 ANO"
$ws.Range("F1").Value = "%SKU COM HISTORICO < 1 ANO"

$ws.Range("O1").Value = "TOTAL SKU INATIVO (ESTOQUE > 0)"
$ws.Range("P1").Value = "%SKU INATIVO (ESTOQUE > 0)"
$ws.Range("Q1").Value = "CUSTO TOTAL INATIVO (ESTOQUE > 0)"
$ws.Range("R1").Value = "TOTAL SKU INATIVO (ESTOQUE <= 0)"
$ws.Range("S1").Value = "%SKU INATIVO (ESTOQUE <= 0)"
$ws.Range("T1").Value = "TOTAL SKU ATIVO (ESTOQUE > 0)"
$ws.Range("U1").Value = "%SKU ATIVO (ESTOQUE > 0)"
$ws.Range("V1").Value = "CUSTO TOTAL ATIVO (ESTOQUE > 0)"
$ws.Range("W1").Value = "TOTAL SKU ATIVO (ESTOQUE <= 0)"
$ws.Range("X1").Value = "%SKU ATIVO (ESTOQUE <= 0)"
$ws.Range("Y1").Value = "TOTAL SKU NAO COMERCIALIZADO (ESTOQUE > 0)"
$ws.Range("Z1").Value = "%SKU NAO COMERCIALIZADO (ESTOQUE > 0)"
$ws.Range("AA1").Value = "CUSTO TOTAL NAO COMERCIALIZADO (ESTOQUE > 0)"
$ws.Range("AB1").Value = "TOTAL SKU NAO COMERCIALIZADO (ESTOQUE <= 0)"
$ws.Range("AC1").Value = "%SKU NAO COMERCIALIZADO (ESTOQUE <= 0)"
$ws.Range("AD1").Value = "TOTAL SKU VERIFICADOS"
$ws.Range("AE1").Value = "TOTAL SKU CONSISTENTES"
$ws.Range("AF1").Value = "%SKU CONSISTENTES"
$ws.Range("AG1").Value = "TOTAL SKU INCONSISTENTES"
$ws.Range("AH1").Value = "%SKU INCONSISTENTES"

# --- Row 2: updated timestamp and metric values ---------------------------
$ws.Range("A2").Value = "2025-05-14 14:32:31"

$ws.Range("C2").Value = 11664
$ws.Range("D2").Value = 72.218438486781
$ws.Range("E2").Value = 2277
$ws.Range("F2").Value = 14.09819825397808

$ws.Range("R2").Value = 3520
$ws.Range("S2").Value = 21.79431614141539

$ws.Range("W2").Value = 2305
$ws.Range("X2").Value = 14.2715621323757

$ws.Range("AI2").Value = 456
$ws.Range("AJ2").Value = 869
$ws.Range("AK2").Value = 1390
$ws.Range("AL2").Value = 16.79558011049724
$ws.Range("AM2").Value = 32.00736648250461
$ws.Range("AN2").Value = 51.19705340699816
$ws.Range("AO2").Value = 1421923.25
$ws.Range("AP2").Value = 266876.46
$ws.Range("AQ2").Value = 88904.47
$ws.Range("AR2").Value = 79.98649415337484
$ws.Range("AS2").Value = 15.01242237052061
$ws.Range("AT2").Value = 5.001083476104556
$ws.Range("AU2").Value = 44.79682953747366
$ws.Range("AV2").Value = 207.5719717544813
$ws.Range("AW2").Value = 567.1548998946259
